# Update the "想去人数" (want-to-go count) values in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        2  = 5672
        3  = 7678
        4  = 80
        7  = 71
        8  = 43
        9  = 4459
        10 = 1797
        13 = 3006
        15 = 574
        16 = 226
        17 = 560
        18 = 484
        19 = 487
        20 = 345
        21 = 123
        22 = 1733
        23 = 1264
        24 = 106
        25 = 1470
        26 = 119
        27 = 595
        31 = 26
        32 = 72
        33 = 112
        34 = 76
        35 = 3205
        36 = 722
        37 = 49
        38 = 185
        40 = 1212
    }
    "演出" = @{
        2 = 15
        3 = 27
    }
    "全部类型" = @{
        2  = 5672
        3  = 7678
        4  = 80
        7  = 71
        8  = 43
        9  = 4459
        10 = 1797
        13 = 3006
        15 = 574
        16 = 226
        17 = 560
        18 = 484
        19 = 487
        20 = 15
        21 = 345
        22 = 123
        23 = 1733
        24 = 1264
        25 = 106
        26 = 1470
        27 = 119
        28 = 595
        32 = 26
        33 = 72
        34 = 112
        35 = 76
        36 = 3205
        37 = 27
        38 = 722
        39 = 49
        40 = 185
        42 = 1212
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
